$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.36078833333333
$ws.Range("N2").Value = 103.082365
$ws.Range("O2").Value = 0.28490270239021
$ws.Range("P2").Value = 0.28490270239021
$ws.Range("Q2").Value = 34.54095340016111
$ws.Range("R2").Value = 310.8685806014499
$ws.Range("S2").Value = 0.02058386203526802
$ws.Range("T2").Value = 0.02058386203526801

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 21.54461566666667
$ws.Range("N3").Value = 64.633847
$ws.Range("O3").Value = 0.17863732245739
$ws.Range("P3").Value = 0.1786373224573899
$ws.Range("Q3").Value = 21.65758126814556
$ws.Range("R3").Value = 194.91823141331
$ws.Range("S3").Value = 0.01290632194417175
$ws.Range("T3").Value = 0.01290632194417174

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.03138866666666
$ws.Range("N4").Value = 180.094166
$ws.Range("O4").Value = 0.4977506538398792
$ws.Range("P4").Value = 0.4977506538398792
$ws.Range("Q4").Value = 60.34615324790888
$ws.Range("R4").Value = 543.1153792311799
$ws.Range("S4").Value = 0.03596185891059694
$ws.Range("T4").Value = 0.03596185891059694

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.005243333333333
$ws.Range("H5").Value = 3.01573
$ws.Range("I5").Value = 0.07224874268505826
$ws.Range("J5").Value = 0.07224874268505825
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.668551
$ws.Range("N5").Value = 14.005653
$ws.Range("O5").Value = 0.03870932131252084
$ws.Range("P5").Value = 0.03870932131252084
$ws.Range("Q5").Value = 4.693029769076666
$ws.Range("R5").Value = 42.23726792169
$ws.Range("S5").Value = 0.00279669979502156
$ws.Range("T5").Value = 0.00279669979502156

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 34.36078833333333
$ws.Range("N6").Value = 103.082365
$ws.Range("O6").Value = 0.28490270239021
$ws.Range("P6").Value = 0.28490270239021
$ws.Range("Q6").Value = 352.5361676666744
$ws.Range("R6").Value = 3172.82550900007
$ws.Range("S6").Value = 0.2100855686762571
$ws.Range("T6").Value = 0.210085568676257

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.54461566666667
$ws.Range("N7").Value = 64.633847
$ws.Range("O7").Value = 0.17863732245739
$ws.Range("P7").Value = 0.1786373224573899
$ws.Range("Q7").Value = 221.0442952384163
$ws.Range("R7").Value = 1989.398657145746
$ws.Range("S7").Value = 0.1317261056508472
$ws.Range("T7").Value = 0.1317261056508471

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7373940889775011
$ws.Range("J8").Value = 0.737394088977501
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.03138866666666
$ws.Range("N8").Value = 180.094166
$ws.Range("O8").Value = 0.4977506538398792
$ws.Range("P8").Value = 0.4977506538398792
$ws.Range("Q8").Value = 615.9124026768876
$ws.Range("R8").Value = 5543.211624091988
$ws.Range("S8").Value = 0.3670383899262132
$ws.Range("T8").Value = 0.3670383899262132

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7373940889775011
$ws.Range("J9").Value = 0.737394088977501
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.668551
$ws.Range("N9").Value = 14.005653
$ws.Range("O9").Value = 0.03870932131252084
$ws.Range("P9").Value = 0.03870932131252084
$ws.Range("Q9").Value = 47.89858317947267
$ws.Range("R9").Value = 431.087248615254
$ws.Range("S9").Value = 0.02854402472418367
$ws.Range("T9").Value = 0.02854402472418367

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 34.36078833333333
$ws.Range("N10").Value = 103.082365
$ws.Range("O10").Value = 0.28490270239021
$ws.Range("P10").Value = 0.28490270239021
$ws.Range("Q10").Value = 13.19024762145833
$ws.Range("R10").Value = 118.712228593125
$ws.Range("S10").Value = 0.007860415261434391
$ws.Range("T10").Value = 0.007860415261434389

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.383875
$ws.Range("H11").Value = 1.151625
$ws.Range("I11").Value = 0.02758982345723265
$ws.Range("J11").Value = 0.02758982345723265
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 21.54461566666667
$ws.Range("N11").Value = 64.633847
$ws.Range("O11").Value = 0.17863732245739
$ws.Range("P11").Value = 0.1786373224573899
$ws.Range("Q11").Value = 8.270439339041667
$ws.Range("R11").Value = 74.43395405137501
$ws.Range("S11").Value = 0.004928572189472132
$ws.Range("T11").Value = 0.004928572189472131

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.383875
$ws.Range("H12").Value = 1.151625
$ws.Range("I12").Value = 0.02758982345723265
$ws.Range("J12").Value = 0.02758982345723265
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.03138866666666
$ws.Range("N12").Value = 180.094166
$ws.Range("O12").Value = 0.4977506538398792
$ws.Range("P12").Value = 0.4977506538398792
$ws.Range("Q12").Value = 23.04454932441667
$ws.Range("R12").Value = 207.40094391975
$ws.Range("S12").Value = 0.01373285266516439
$ws.Range("T12").Value = 0.01373285266516439

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.383875
$ws.Range("H13").Value = 1.151625
$ws.Range("I13").Value = 0.02758982345723265
$ws.Range("J13").Value = 0.02758982345723265
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.668551
$ws.Range("N13").Value = 14.005653
$ws.Range("O13").Value = 0.03870932131252084
$ws.Range("P13").Value = 0.03870932131252084
$ws.Range("Q13").Value = 1.792140015125
$ws.Range("R13").Value = 16.129260136125
$ws.Range("S13").Value = 0.001067983341161743
$ws.Range("T13").Value = 0.001067983341161743

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.264687
$ws.Range("H14").Value = 6.794061
$ws.Range("I14").Value = 0.162767344880208
$ws.Range("J14").Value = 0.162767344880208
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 34.36078833333333
$ws.Range("N14").Value = 103.082365
$ws.Range("O14").Value = 0.28490270239021
$ws.Range("P14").Value = 0.28490270239021
$ws.Range("Q14").Value = 77.81643064825165
$ws.Range("R14").Value = 700.347875834265
$ws.Range("S14").Value = 0.04637285641725058
$ws.Range("T14").Value = 0.04637285641725056

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.264687
$ws.Range("H15").Value = 6.794061
$ws.Range("I15").Value = 0.162767344880208
$ws.Range("J15").Value = 0.162767344880208
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 21.54461566666667
$ws.Range("N15").Value = 64.633847
$ws.Range("O15").Value = 0.17863732245739
$ws.Range("P15").Value = 0.1786373224573899
$ws.Range("Q15").Value = 48.79181102029634
$ws.Range("R15").Value = 439.126299182667
$ws.Range("S15").Value = 0.02907632267289892
$ws.Range("T15").Value = 0.02907632267289892

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.264687
$ws.Range("H16").Value = 6.794061
$ws.Range("I16").Value = 0.162767344880208
$ws.Range("J16").Value = 0.162767344880208
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.03138866666666
$ws.Range("N16").Value = 180.094166
$ws.Range("O16").Value = 0.4977506538398792
$ws.Range("P16").Value = 0.4977506538398792
$ws.Range("Q16").Value = 135.9523055053473
$ws.Range("R16").Value = 1223.570749548126
$ws.Range("S16").Value = 0.08101755233790464
$ws.Range("T16").Value = 0.08101755233790464

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.264687
$ws.Range("H17").Value = 6.794061
$ws.Range("I17").Value = 0.162767344880208
$ws.Range("J17").Value = 0.162767344880208
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.668551
$ws.Range("N17").Value = 14.005653
$ws.Range("O17").Value = 0.03870932131252084
$ws.Range("P17").Value = 0.03870932131252084
$ws.Range("Q17").Value = 10.572806758537
$ws.Range("R17").Value = 95.155260826833
$ws.Range("S17").Value = 0.006300613452153866
$ws.Range("T17").Value = 0.006300613452153866
